$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add column B header first so "slug" gets the next shared-string slot
$ws.Range("B1").Value = "slug"

# Update column A values (name, Camera, Recorder), clear A4
$ws.Range("A2").Value = "Camera"
$ws.Range("B2").Value = "camera"
$ws.Range("A3").Value = "Recorder"
$ws.Range("B3").Value = "recorder"
$ws.Range("A4").Value = ""

# Update selection to C2
$ws.Range("C2").Select()
